# Week 13 logging update for Bears Players Data
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# Row 2: D.Montgomery
$rushing.Range("D2").Value = 2
$rushing.Range("F2").Value = 2

# Row 3: D.Williams
$rushing.Range("C3").Value = 111
$rushing.Range("D3").Value = 53
$rushing.Range("E3").Value = 23
$rushing.Range("F3").Value = 19

# Row 5: R.Nall
$rushing.Range("C5").Value = 15
$rushing.Range("D5").Value = 7
$rushing.Range("E5").Value = 2

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# Row 2: D.Williams
$receiving.Range("C2").Value = 33
$receiving.Range("D2").Value = 29
$receiving.Range("G2").Value = 5
$receiving.Range("H2").Value = 5

# Row 4: M.Goodwin
$receiving.Range("C4").Value = 4
$receiving.Range("D4").Value = 4

# Row 6: D.Mooney
$receiving.Range("C6").Value = 67
$receiving.Range("D6").Value = 39
$receiving.Range("G6").Value = 5
$receiving.Range("H6").Value = 3

# Row 8: D.Byrd
$receiving.Range("C8").Value = 16
$receiving.Range("D8").Value = 12
$receiving.Range("E8").Value = 3
$receiving.Range("G8").Value = 2
$receiving.Range("H8").Value = 2

# Row 9: C.Kmet
$receiving.Range("C9").Value = 10
$receiving.Range("D9").Value = 7
$receiving.Range("E9").Value = 2
$receiving.Range("F9").Value = 1
$receiving.Range("G9").Value = 1
$receiving.Range("H9").Value = 1

# Row 10: J.Graham
$receiving.Range("C10").Value = 55
$receiving.Range("D10").Value = 35
$receiving.Range("E10").Value = 9
$receiving.Range("F10").Value = 5
$receiving.Range("G10").Value = 8

# Row 11: J.James
$receiving.Range("C11").Value = 11
$receiving.Range("D11").Value = 5
$receiving.Range("G11").Value = 6
$receiving.Range("H11").Value = 3
